# "Préparation du cours 4"
# The slide "RÉGLER UNE TENSION AVEC digitalWrite()" (first slide, sldId 2941)
# is no longer needed and is removed from the deck, leaving only the
# remaining slide (sldId 3070) in the presentation.

$p = $ppt.ActivePresentation

$p.Slides.Item(1).Delete()
